$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Employment (% of total)" row (row 12): Micro / SMEs / MSMEs values.
# These are stored as literal text (not numbers) in the workbook, so we
# briefly force a Text number format while assigning the value (otherwise
# Excel auto-converts a numeric-looking string to a real number), then
# restore the original "Normal" style so no stray formatting is left behind.

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "78.35"
$ws.Range("B12").Style = "Normal"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "20.24"
$ws.Range("C12").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.59"
$ws.Range("D12").Style = "Normal"
